# Update "想去人数" (F column) counters across the four sheets to match the
# freshly generated gh-pages data snapshot (commit 456a3b4).

$wb = $excel.ActiveWorkbook

# --- 展览 (Exhibitions) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Cells.Item(2, 6).Value = 56
$ws1.Cells.Item(4, 6).Value = 5923
$ws1.Cells.Item(6, 6).Value = 2978
$ws1.Cells.Item(8, 6).Value = 403
$ws1.Cells.Item(12, 6).Value = 713
$ws1.Cells.Item(13, 6).Value = 248
$ws1.Cells.Item(14, 6).Value = 4333
$ws1.Cells.Item(15, 6).Value = 4333
$ws1.Cells.Item(18, 6).Value = 112
$ws1.Cells.Item(22, 6).Value = 6623
$ws1.Cells.Item(23, 6).Value = 229
$ws1.Cells.Item(24, 6).Value = 100
$ws1.Cells.Item(25, 6).Value = 293
$ws1.Cells.Item(26, 6).Value = 456
$ws1.Cells.Item(27, 6).Value = 1240
$ws1.Cells.Item(28, 6).Value = 6249
$ws1.Cells.Item(29, 6).Value = 1635
$ws1.Cells.Item(30, 6).Value = 14
$ws1.Cells.Item(31, 6).Value = 1863
$ws1.Cells.Item(32, 6).Value = 5981
$ws1.Cells.Item(33, 6).Value = 112
$ws1.Cells.Item(35, 6).Value = 94
$ws1.Cells.Item(37, 6).Value = 416
$ws1.Cells.Item(38, 6).Value = 4118
$ws1.Cells.Item(39, 6).Value = 13
$ws1.Cells.Item(40, 6).Value = 187
$ws1.Cells.Item(41, 6).Value = 82
$ws1.Cells.Item(43, 6).Value = 2405
$ws1.Cells.Item(44, 6).Value = 24
$ws1.Cells.Item(47, 6).Value = 18
$ws1.Cells.Item(48, 6).Value = 327
$ws1.Cells.Item(49, 6).Value = 2052
$ws1.Cells.Item(50, 6).Value = 23

# --- 演出 (Performances) ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Cells.Item(6, 6).Value = 103

# --- 本地生活 (Local life) ---
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Cells.Item(2, 6).Value = 1414

# --- 全部类型 (All types) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Cells.Item(2, 6).Value = 56
$ws4.Cells.Item(3, 6).Value = 1414
$ws4.Cells.Item(5, 6).Value = 5923
$ws4.Cells.Item(7, 6).Value = 2978
$ws4.Cells.Item(14, 6).Value = 248
$ws4.Cells.Item(15, 6).Value = 4333
$ws4.Cells.Item(16, 6).Value = 4333
$ws4.Cells.Item(19, 6).Value = 112
$ws4.Cells.Item(22, 6).Value = 6623
$ws4.Cells.Item(23, 6).Value = 229
$ws4.Cells.Item(24, 6).Value = 100
$ws4.Cells.Item(25, 6).Value = 456
$ws4.Cells.Item(26, 6).Value = 1241
$ws4.Cells.Item(27, 6).Value = 103
$ws4.Cells.Item(28, 6).Value = 6249
$ws4.Cells.Item(29, 6).Value = 1635
$ws4.Cells.Item(31, 6).Value = 14
$ws4.Cells.Item(32, 6).Value = 1863
$ws4.Cells.Item(33, 6).Value = 5981
$ws4.Cells.Item(34, 6).Value = 112
$ws4.Cells.Item(36, 6).Value = 94
$ws4.Cells.Item(38, 6).Value = 416
$ws4.Cells.Item(39, 6).Value = 4118
$ws4.Cells.Item(40, 6).Value = 13
$ws4.Cells.Item(41, 6).Value = 187
$ws4.Cells.Item(42, 6).Value = 82
$ws4.Cells.Item(45, 6).Value = 2405
$ws4.Cells.Item(46, 6).Value = 24
$ws4.Cells.Item(49, 6).Value = 18

$wb.Save()
